$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Objetivos -> value becomes the "Danubia" text (content shift)
$ws.Range("B10").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C10").Value = "9146830 - Danúbia Caporusso Bargos"

# Row 12 stays "Docentes responsáveis:" with no B/C (unchanged)

# Row 13 (was "9146830..." row with no A) becomes "Programa resumido:" row,
# with B/C now holding "01/01/2020" (reused text) and row height 60
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2020"
$ws.Range("C13").Value = "01/01/2020"
$ws.Rows.Item(13).RowHeight = 45

# Row 14 (was "5464150..." row with no A) becomes "Short syllabus:" row
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Introduction to cartography. Basic concepts of cartography and geodesy. Scales. Coordinate systems. Cartographic projections. National Cartographic System (NCS). Introduction to Remote Sensing.  Global Positioning System (GPS). Digital cartography. Use of charts, maps and plants in environmental studies. Introduction to topography. Methods to obtain topographic measurements. Topographic survey. Leveling. Topographic maps in environmental studies."
$ws.Range("C14").Value = $ws.Range("B14").Value
$ws.Rows.Item(14).RowHeight = 45

# Row 15 (was "Programa resumido:") becomes "Programa:" with reused Danubia text
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C15").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Rows.Item(15).RowHeight = 90

# Row 16 (was "Short syllabus:") becomes "Syllabus:"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Introduction to cartography and the history of maps; definitions and basic concepts of cartography and geodesy; shape and dimensions of land; coordinate systems. Cartographic projections. Scale and cartographic accuracy. Cartographic representation. Thematic cartography. Introduction to Remote Sensing. ""Global Navigation Satellite System"" - GNSS. Introduction to topography: reference surfaces, errors of observation. Technical standard NBR 13.133. Methods to obtain topographic measurements. Measurements of distances: direct and indirect measures. Direction measurements: horizontal and vertical angles. Orientation: magnetic north; geographic north; azimuth; course. Topographic survey - planimetry. Leveling. Topographic maps applied to environmental studies. Representation of landforms. Maps applied to environmental studies. Spatial data infrastructure. Digital cartography."
$ws.Range("C16").Value = $ws.Range("B16").Value
$ws.Rows.Item(16).RowHeight = 90

# Row 17 (was "Programa:") becomes "Avaliação:" with B/C cleared
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Rows.Item(17).RowHeight = 15

# Row 18 (was "Syllabus:") becomes "Método:" with Mariana text
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C18").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Rows.Item(18).RowHeight = 45

# Row 19 (was "Avaliação:" no B/C) becomes "Critério:" with "Aulas expositivas..." text
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas teóricas e práticas, visitas técnicas, pesquisas bibliográficas e estudos dirigidos."
$ws.Range("C19").Value = $ws.Range("B19").Value
$ws.Rows.Item(19).RowHeight = 45

# Row 20 (was "Método:" with "Aulas expositivas..." text) becomes "Norma de recuperação:"
# with "Média ponderada..." text
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada de exercícios e provas."
$ws.Range("C20").Value = "Média ponderada de exercícios e provas."
$ws.Rows.Item(20).RowHeight = 45

# Row 21 (was "Critério:" with "Média ponderada..." text) becomes "Bibliografia:"
# with the "A nota final..." text, and row height grows to 120 (90pt)
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A nota final (MF) do aluno que realizou provas de recuperação dependerá da média do semestre (MS) e da média das provas de recuperação (MR), como segue:MF=5 se 5 ≤MR ≤ (10 - MS); MF = (MS + MR) / 2 se MR > (10 – MS)MF = MS se MR < 5."
$ws.Range("C21").Value = $ws.Range("B21").Value
$ws.Rows.Item(21).RowHeight = 90

# Delete old rows 22 and 23 (Norma de recuperação / Bibliografia with the long
# bibliography text), which are no longer present in the sheet.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
